# Added upload api endpoint info.
#
# Order of operations matters here: new shared strings are interned in the
# order cells are first assigned a (new) text value, and the target workbook
# needs a very specific shared-string ordering, so values are written in the
# sequence below rather than strictly row-by-row / left-to-right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- /user/login row (row 4): "Body" cell becomes a two-line value ---
# Was a flat "Email,Password" -> becomes wrapped "Email" / "Password" and the
# row grows taller to fit the wrapped text.
$ws.Range("E4").Value = "Email" + [char]10 + "Password"
$ws.Range("E4").WrapText = $true
$ws.Rows.Item(4).RowHeight = 30

# --- New row 6: POST /upload/file ---
$ws.Range("C6").Value = "add single or multiple files to the server"
$ws.Range("D6").Value = "yes"
$ws.Range("F6").Value = "The path to the file "
$ws.Range("E6").Value = "file(single/multiple)"

# --- New row 7: GET /upload/ ---
$ws.Range("A7").Value = "/upload/"
$ws.Range("A6").Value = "/upload/file"
$ws.Range("C7").Value = "get a list of all the files"
$ws.Range("F7").Value = "All the files in the uploads folder"

$ws.Range("B6").Value = "post"
$ws.Range("B7").Value = "get"
$ws.Range("D7").Value = "yes"

# --- Columns E & F grew wider to fit the new content ---
$ws.Columns.Item(5).ColumnWidth = 18.6
$ws.Columns.Item(6).ColumnWidth = 29.8

# --- Selection ends on the newly added E7 cell ---
$ws.Range("E7").Select()
